$wb = $excel.ActiveWorkbook

# --- Results sheet updates ---
$ws = $wb.Worksheets.Item("Results")

# Ensure percentage-looking text values (columns J and S) stay as text, not
# auto-converted by Excel's percentage number-format detection.
$ws.Range("J2:J4").NumberFormat = "@"
$ws.Range("J6").NumberFormat = "@"
$ws.Range("S2:S4").NumberFormat = "@"
$ws.Range("S6").NumberFormat = "@"

# Row 2 (Bid ID 1)
$ws.Range("F2").Value = 1
$ws.Range("H2").Value = "Supplier 1"
$ws.Range("I2").Value = 34
$ws.Range("J2").Value = "2%"
$ws.Range("K2").Value = 33.32
$ws.Range("L2").Value = "DDP"
$ws.Range("M2").Value = 2
$ws.Range("N2").Value = 35.32
$ws.Range("O2").Value = 51425.92
$ws.Range("Q2").Value = -15025.92
$ws.Range("R2").Value = -49969.92
$ws.Range("S2").Value = "95%"
$ws.Range("T2").Value = 48854.624

# Row 3 (Bid ID 2)
$ws.Range("F3").Value = 1
$ws.Range("H3").Value = "Supplier 1"
$ws.Range("I3").Value = 15
$ws.Range("J3").Value = "2%"
$ws.Range("K3").Value = 14.7
$ws.Range("L3").Value = "DDP"
$ws.Range("M3").Value = 2
$ws.Range("N3").Value = 16.7
$ws.Range("O3").Value = 40430.7
$ws.Range("Q3").Value = 46725.3
$ws.Range("R3").Value = -38009.7
$ws.Range("S3").Value = "95%"
$ws.Range("T3").Value = 38409.16499999999

# Row 4 (Bid ID 3)
$ws.Range("F4").Value = 1
$ws.Range("H4").Value = "Supplier 1"
$ws.Range("I4").Value = 15
$ws.Range("J4").Value = "2%"
$ws.Range("K4").Value = 14.7
$ws.Range("L4").Value = "DDP"
$ws.Range("M4").Value = 2
$ws.Range("N4").Value = 16.7
$ws.Range("O4").Value = 107915.4
$ws.Range("Q4").Value = 421968.6
$ws.Range("R4").Value = -101453.4
$ws.Range("S4").Value = "95%"
$ws.Range("T4").Value = 102519.63

# Row 5 (Bid ID 4) - only Current Price and Current Price Savings change
$ws.Range("F5").Value = 1
$ws.Range("R5").Value = -175596.5

# Row 6 (Bid ID 5)
$ws.Range("F6").Value = 1
$ws.Range("H6").Value = "Supplier 1"
$ws.Range("I6").Value = 24
$ws.Range("J6").Value = "2%"
$ws.Range("K6").Value = 23.52
$ws.Range("L6").Value = "DDP"
$ws.Range("M6").Value = 2
$ws.Range("N6").Value = 25.52
$ws.Range("O6").Value = 61783.92
$ws.Range("Q6").Value = 119791.08
$ws.Range("R6").Value = -59362.92
$ws.Range("S6").Value = "95%"
$ws.Range("T6").Value = 58694.72399999999

# --- LP Model sheet update ---
$wsLp = $wb.Worksheets.Item("LP Model")
$wsLp.Range("A2").Value = "\* Sourcing_with_Freight *\`nMinimize`nOBJ: F_Supplier_1 + F_Supplier_2 + F_Supplier_3 + F_Supplier_4 + F_Supplier_5`n + S_Supplier_1 + S_Supplier_2 + S_Supplier_3 + S_Supplier_4 + S_Supplier_5`n - reb_Supplier_1 - reb_Supplier_2 - reb_Supplier_3 - reb_Supplier_4`n - reb_Supplier_5`nSubject To`nCap_Supplier_1_Bid_ID_1: x_Supplier_1_1 <= 1457`nCap_Supplier_1_Bid_ID_2: x_Supplier_1_2 <= 2422`nCap_Supplier_1_Bid_ID_3: x_Supplier_1_3 <= 6463`nCap_Supplier_1_Bid_ID_4: x_Supplier_1_4 <= 2358`nCap_Supplier_1_Bid_ID_5: x_Supplier_1_5 <= 2422`nCap_Supplier_2_Capacity_Group_Category_1: x_Supplier_2_1 + x_Supplier_2_3`n + x_Supplier_2_4 <= 8001`nCap_Supplier_2_Capacity_Group_Category_2: x_Supplier_2_2 <= 2501`nCap_Supplier_2_Capacity_Group_Category_3: x_Supplier_2_5 <= 501`nCap_Supplier_3_Bid_ID_1: x_Supplier_3_1 <= 1457`nCap_Supplier_3_Bid_ID_2: x_Supplier_3_2 <= 2422`nCap_Supplier_3_Bid_ID_3: x_Supplier_3_3 <= 6463`nCap_Supplier_3_Bid_ID_4: x_Supplier_3_4 <= 2358`nCap_Supplier_3_Bid_ID_5: x_Supplier_3_5 <= 2422`nCap_Supplier_4_Bid_ID_1: x_Supplier_4_1 <= 1457`nCap_Supplier_4_Bid_ID_2: x_Supplier_4_2 <= 2422`nCap_Supplier_4_Bid_ID_3: x_Supplier_4_3 <= 6463`nCap_Supplier_4_Bid_ID_4: x_Supplier_4_4 <= 2358`nCap_Supplier_4_Bid_ID_5: x_Supplier_4_5 <= 2422`nCap_Supplier_5_Description_Large_Item: x_Supplier_5_5 <= 5001`nCap_Supplier_5_Description_Medium_item: x_Supplier_5_2 <= 2501`nCap_Supplier_5_Description_Small_item: x_Supplier_5_1 + x_Supplier_5_3`n + x_Supplier_5_4 <= 9001`nDDP_ON_Supplier_1_1: - 1000000000 bDDP_Supplier_1_1 + xSF_Supplier_1_1 <= 0`nDDP_ON_Supplier_1_2: - 1000000000 bDDP_Supplier_1_2 + xSF_Supplier_1_2 <= 0`nDDP_ON_Supplier_1_3: - 1000000000 bDDP_Supplier_1_3 + xSF_Supplier_1_3 <= 0`nDDP_ON_Supplier_1_4: - 1000000000 bDDP_Supplier_1_4 + xSF_Supplier_1_4 <= 0`nDDP_ON_Supplier_1_5: - 1000000000 bDDP_Supplier_1_5 + xSF_Supplier_1_5 <= 0`nDDP_ON_Supplier_2_1: - 1000000000 bDDP_Supplier_2_1 + xSF_Supplier_2_1 <= 0`nDDP_ON_Supplier_2_2: - 1000000000 bDDP_Supplier_2_2 + xSF_Supplier_2_2 <= 0`nDDP_ON_Supplier_2_3: - 1000000000 bDDP_Supplier_2_3 + xSF_Supplier_2_3 <= 0`nDDP_ON_Supplier_2_4: - 1000000000 bDDP_Supplier_2_4 + xSF_Supplier_2_4 <= 0`nDDP_ON_Supplier_2_5: - 1000000000 bDDP_Supplier_2_5 + xSF_Supplier_2_5 <= 0`nDDP_ON_Supplier_3_1: - 1000000000 bDDP_Supplier_3_1 + xSF_Supplier_3_1 <= 0`nDDP_ON_Supplier_3_2: - 1000000000 bDDP_Supplier_3_2 + xSF_Supplier_3_2 <= 0`nDDP_ON_Supplier_3_3: - 1000000000 bDDP_Supplier_3_3 + xSF_Supplier_3_3 <= 0`nDDP_ON_Supplier_3_4: - 1000000000 bDDP_Supplier_3_4 + xSF_Supplier_3_4 <= 0`nDDP_ON_Supplier_3_5: - 1000000000 bDDP_Supplier_3_5 + xSF_Supplier_3_5 <= 0`nDDP_ON_Supplier_4_1: - 1000000000 bDDP_Supplier_4_1 + xSF_Supplier_4_1 <= 0`nDDP_ON_Supplier_4_2: - 1000000000 bDDP_Supplier_4_2 + xSF_Supplier_4_2 <= 0`nDDP_ON_Supplier_4_3: - 1000000000 bDDP_Supplier_4_3 + xSF_Supplier_4_3 <= 0`nDDP_ON_Supplier_4_4: - 1000000000 bDDP_Supplier_4_4 + xSF_Supplier_4_4 <= 0`nDDP_ON_Supplier_4_5: - 1000000000 bDDP_Supplier_4_5 + xSF_Supplier_4_5 <= 0`nDDP_ON_Supplier_5_1: - 1000000000 bDDP_Supplier_5_1 + xSF_Supplier_5_1 <= 0`nDDP_ON_Supplier_5_2: - 1000000000 bDDP_Supplier_5_2 + xSF_Supplier_5_2 <= 0`nDDP_ON_Supplier_5_3: - 1000000000 bDDP_Supplier_5_3 + xSF_Supplier_5_3 <= 0`nDDP_ON_Supplier_5_4: - 1000000000 bDDP_Supplier_5_4 + xSF_Supplier_5_4 <= 0`nDDP_ON_Supplier_5_5: - 1000000000 bDDP_Supplier_5_5 + xSF_Supplier_5_5 <= 0`nDMax_Supplier_1_0: x_Supplier_1_1 + x_Supplier_1_2 + x_Supplier_1_3`n + x_Supplier_1_4 + x_Supplier_1_5 + 1000000000 zd_Supplier_1_0 <= 1001000000`nDMin_Supplier_1_0: x_Supplier_1_1 + x_Supplier_1_2 + x_Supplier_1_3`n + x_Supplier_1_4 + x_Supplier_1_5 - zd_Supplier_1_0 >= 0`nDemand_1: x_Supplier_1_1 + x_Supplier_2_1 + x_Supplier_3_1 + x_Supplier_4_1`n + x_Supplier_5_1 = 1456`nDemand_2: x_Supplier_1_2 + x_Supplier_2_2 + x_Supplier_3_2 + x_Supplier_4_2`n + x_Supplier_5_2 = 2421`nDemand_3: x_Supplier_1_3 + x_Supplier_2_3 + x_Supplier_3_3 + x_Supplier_4_3`n + x_Supplier_5_3 = 6462`nDemand_4: x_Supplier_1_4 + x_Supplier_2_4 + x_Supplier_3_4 + x_Supplier_4_4`n + x_Supplier_5_4 = 2357`nDemand_5: x_Supplier_1_5 + x_Supplier_2_5 + x_Supplier_3_5 + x_Supplier_4_5`n + x_Supplier_5_5 = 2421`nF_Supplier_1: F_Supplier_1 - xKBX_Supplier_1_1 - xKBX_Supplier_1_2`n - xKBX_Supplier_1_3 - xKBX_Supplier_1_4 - xKBX_Supplier_1_5 = 0`nF_Supplier_2: F_Supplier_2 - xKBX_Supplier_2_1 - xKBX_Supplier_2_2`n - xKBX_Supplier_2_3 - xKBX_Supplier_2_4 - xKBX_Supplier_2_5 = 0`nF_Supplier_3: F_Supplier_3 - xKBX_Supplier_3_1 - xKBX_Supplier_3_2`n - xKBX_Supplier_3_3 - xKBX_Supplier_3_4 - xKBX_Supplier_3_5 = 0`nF_Supplier_4: F_Supplier_4 - xKBX_Supplier_4_1 - xKBX_Supplier_4_2`n - xKBX_Supplier_4_3 - xKBX_Supplier_4_4 - xKBX_Supplier_4_5 = 0`nF_Supplier_5: F_Supplier_5 - xKBX_Supplier_5_1 - xKBX_Supplier_5_2`n - xKBX_Supplier_5_3 - xKBX_Supplier_5_4 - xKBX_Supplier_5_5 = 0`nFixd_Supplier_2: d_Supplier_2 = 0`nFixd_Supplier_3: d_Supplier_3 = 0`nFixd_Supplier_4: d_Supplier_4 = 0`nFixd_Supplier_5: d_Supplier_5 = 0`nFixreb_Supplier_2: reb_Supplier_2 = 0`nFixreb_Supplier_3: reb_Supplier_3 = 0`nFixreb_Supplier_4: reb_Supplier_4 = 0`nFixreb_Supplier_5: reb_Supplier_5 = 0`nKBX_ON_Supplier_1_1: 1000000000 bDDP_Supplier_1_1 + xKBX_Supplier_1_1`n <= 1000000000`nKBX_ON_Supplier_1_2: 1000000000 bDDP_Supplier_1_2 + xKBX_Supplier_1_2`n <= 1000000000`nKBX_ON_Supplier_1_3: 1000000000 bDDP_Supplier_1_3 + xKBX_Supplier_1_3`n <= 1000000000`nKBX_ON_Supplier_1_4: 1000000000 bDDP_Supplier_1_4 + xKBX_Supplier_1_4`n <= 1000000000`nKBX_ON_Supplier_1_5: 1000000000 bDDP_Supplier_1_5 + xKBX_Supplier_1_5`n <= 1000000000`nKBX_ON_Supplier_2_1: 1000000000 bDDP_Supplier_2_1 + xKBX_Supplier_2_1`n <= 1000000000`nKBX_ON_Supplier_2_2: 1000000000 bDDP_Supplier_2_2 + xKBX_Supplier_2_2`n <= 1000000000`nKBX_ON_Supplier_2_3: 1000000000 bDDP_Supplier_2_3 + xKBX_Supplier_2_3`n <= 1000000000`nKBX_ON_Supplier_2_4: 1000000000 bDDP_Supplier_2_4 + xKBX_Supplier_2_4`n <= 1000000000`nKBX_ON_Supplier_2_5: 1000000000 bDDP_Supplier_2_5 + xKBX_Supplier_2_5`n <= 1000000000`nKBX_ON_Supplier_3_1: 1000000000 bDDP_Supplier_3_1 + xKBX_Supplier_3_1`n <= 1000000000`nKBX_ON_Supplier_3_2: 1000000000 bDDP_Supplier_3_2 + xKBX_Supplier_3_2`n <= 1000000000`nKBX_ON_Supplier_3_3: 1000000000 bDDP_Supplier_3_3 + xKBX_Supplier_3_3`n <= 1000000000`nKBX_ON_Supplier_3_4: 1000000000 bDDP_Supplier_3_4 + xKBX_Supplier_3_4`n <= 1000000000`nKBX_ON_Supplier_3_5: 1000000000 bDDP_Supplier_3_5 + xKBX_Supplier_3_5`n <= 1000000000`nKBX_ON_Supplier_4_1: 1000000000 bDDP_Supplier_4_1 + xKBX_Supplier_4_1`n <= 1000000000`nKBX_ON_Supplier_4_2: 1000000000 bDDP_Supplier_4_2 + xKBX_Supplier_4_2`n <= 1000000000`nKBX_ON_Supplier_4_3: 1000000000 bDDP_Supplier_4_3 + xKBX_Supplier_4_3`n <= 1000000000`nKBX_ON_Supplier_4_4: 1000000000 bDDP_Supplier_4_4 + xKBX_Supplier_4_4`n <= 1000000000`nKBX_ON_Supplier_4_5: 1000000000 bDDP_Supplier_4_5 + xKBX_Supplier_4_5`n <= 1000000000`nKBX_ON_Supplier_5_1: 1000000000 bDDP_Supplier_5_1 + xKBX_Supplier_5_1`n <= 1000000000`nKBX_ON_Supplier_5_2: 1000000000 bDDP_Supplier_5_2 + xKBX_Supplier_5_2`n <= 1000000000`nKBX_ON_Supplier_5_3: 1000000000 bDDP_Supplier_5_3 + xKBX_Supplier_5_3`n <= 1000000000`nKBX_ON_Supplier_5_4: 1000000000 bDDP_Supplier_5_4 + xKBX_Supplier_5_4`n <= 1000000000`nKBX_ON_Supplier_5_5: 1000000000 bDDP_Supplier_5_5 + xKBX_Supplier_5_5`n <= 1000000000`nLinkVol_Supplier_1_1: - xKBX_Supplier_1_1 - xSF_Supplier_1_1 + x_Supplier_1_1`n = 0`nLinkVol_Supplier_1_2: - xKBX_Supplier_1_2 - xSF_Supplier_1_2 + x_Supplier_1_2`n = 0`nLinkVol_Supplier_1_3: - xKBX_Supplier_1_3 - xSF_Supplier_1_3 + x_Supplier_1_3`n = 0`nLinkVol_Supplier_1_4: - xKBX_Supplier_1_4 - xSF_Supplier_1_4 + x_Supplier_1_4`n = 0`nLinkVol_Supplier_1_5: - xKBX_Supplier_1_5 - xSF_Supplier_1_5 + x_Supplier_1_5`n = 0`nLinkVol_Supplier_2_1: - xKBX_Supplier_2_1 - xSF_Supplier_2_1 + x_Supplier_2_1`n = 0`nLinkVol_Supplier_2_2: - xKBX_Supplier_2_2 - xSF_Supplier_2_2 + x_Supplier_2_2`n = 0`nLinkVol_Supplier_2_3: - xKBX_Supplier_2_3 - xSF_Supplier_2_3 + x_Supplier_2_3`n = 0`nLinkVol_Supplier_2_4: - xKBX_Supplier_2_4 - xSF_Supplier_2_4 + x_Supplier_2_4`n = 0`nLinkVol_Supplier_2_5: - xKBX_Supplier_2_5 - xSF_Supplier_2_5 + x_Supplier_2_5`n = 0`nLinkVol_Supplier_3_1: - xKBX_Supplier_3_1 - xSF_Supplier_3_1 + x_Supplier_3_1`n = 0`nLinkVol_Supplier_3_2: - xKBX_Supplier_3_2 - xSF_Supplier_3_2 + x_Supplier_3_2`n = 0`nLinkVol_Supplier_3_3: - xKBX_Supplier_3_3 - xSF_Supplier_3_3 + x_Supplier_3_3`n = 0`nLinkVol_Supplier_3_4: - xKBX_Supplier_3_4 - xSF_Supplier_3_4 + x_Supplier_3_4`n = 0`nLinkVol_Supplier_3_5: - xKBX_Supplier_3_5 - xSF_Supplier_3_5 + x_Supplier_3_5`n = 0`nLinkVol_Supplier_4_1: - xKBX_Supplier_4_1 - xSF_Supplier_4_1 + x_Supplier_4_1`n = 0`nLinkVol_Supplier_4_2: - xKBX_Supplier_4_2 - xSF_Supplier_4_2 + x_Supplier_4_2`n = 0`nLinkVol_Supplier_4_3: - xKBX_Supplier_4_3 - xSF_Supplier_4_3 + x_Supplier_4_3`n = 0`nLinkVol_Supplier_4_4: - xKBX_Supplier_4_4 - xSF_Supplier_4_4 + x_Supplier_4_4`n = 0`nLinkVol_Supplier_4_5: - xKBX_Supplier_4_5 - xSF_Supplier_4_5 + x_Supplier_4_5`n = 0`nLinkVol_Supplier_5_1: - xKBX_Supplier_5_1 - xSF_Supplier_5_1 + x_Supplier_5_1`n = 0`nLinkVol_Supplier_5_2: - xKBX_Supplier_5_2 - xSF_Supplier_5_2 + x_Supplier_5_2`n = 0`nLinkVol_Supplier_5_3: - xKBX_Supplier_5_3 - xSF_Supplier_5_3 + x_Supplier_5_3`n = 0`nLinkVol_Supplier_5_4: - xKBX_Supplier_5_4 - xSF_Supplier_5_4 + x_Supplier_5_4`n = 0`nLinkVol_Supplier_5_5: - xKBX_Supplier_5_5 - xSF_Supplier_5_5 + x_Supplier_5_5`n = 0`nLink_Supplier_1_1: x_Supplier_1_1 - 1000000000 z_Supplier_1 <= 0`nLink_Supplier_1_2: x_Supplier_1_2 - 1000000000 z_Supplier_1 <= 0`nLink_Supplier_1_3: x_Supplier_1_3 - 1000000000 z_Supplier_1 <= 0`nLink_Supplier_1_4: x_Supplier_1_4 - 1000000000 z_Supplier_1 <= 0`nLink_Supplier_1_5: x_Supplier_1_5 - 1000000000 z_Supplier_1 <= 0`nLink_Supplier_2_1: x_Supplier_2_1 - 1000000000 z_Supplier_2 <= 0`nLink_Supplier_2_2: x_Supplier_2_2 - 1000000000 z_Supplier_2 <= 0`nLink_Supplier_2_3: x_Supplier_2_3 - 1000000000 z_Supplier_2 <= 0`nLink_Supplier_2_4: x_Supplier_2_4 - 1000000000 z_Supplier_2 <= 0`nLink_Supplier_2_5: x_Supplier_2_5 - 1000000000 z_Supplier_2 <= 0`nLink_Supplier_3_1: x_Supplier_3_1 - 1000000000 z_Supplier_3 <= 0`nLink_Supplier_3_2: x_Supplier_3_2 - 1000000000 z_Supplier_3 <= 0`nLink_Supplier_3_3: x_Supplier_3_3 - 1000000000 z_Supplier_3 <= 0`nLink_Supplier_3_4: x_Supplier_3_4 - 1000000000 z_Supplier_3 <= 0`nLink_Supplier_3_5: x_Supplier_3_5 - 1000000000 z_Supplier_3 <= 0`nLink_Supplier_4_1: x_Supplier_4_1 - 1000000000 z_Supplier_4 <= 0`nLink_Supplier_4_2: x_Supplier_4_2 - 1000000000 z_Supplier_4 <= 0`nLink_Supplier_4_3: x_Supplier_4_3 - 1000000000 z_Supplier_4 <= 0`nLink_Supplier_4_4: x_Supplier_4_4 - 1000000000 z_Supplier_4 <= 0`nLink_Supplier_4_5: x_Supplier_4_5 - 1000000000 z_Supplier_4 <= 0`nLink_Supplier_5_1: x_Supplier_5_1 - 1000000000 z_Supplier_5 <= 0`nLink_Supplier_5_2: x_Supplier_5_2 - 1000000000 z_Supplier_5 <= 0`nLink_Supplier_5_3: x_Supplier_5_3 - 1000000000 z_Supplier_5 <= 0`nLink_Supplier_5_4: x_Supplier_5_4 - 1000000000 z_Supplier_5 <= 0`nLink_Supplier_5_5: x_Supplier_5_5 - 1000000000 z_Supplier_5 <= 0`nMinAward_Supplier_1: x_Supplier_1_1 + x_Supplier_1_2 + x_Supplier_1_3`n + x_Supplier_1_4 + x_Supplier_1_5 - z_Supplier_1 >= 0`nMinAward_Supplier_2: x_Supplier_2_1 + x_Supplier_2_2 + x_Supplier_2_3`n + x_Supplier_2_4 + x_Supplier_2_5 - z_Supplier_2 >= 0`nMinAward_Supplier_3: x_Supplier_3_1 + x_Supplier_3_2 + x_Supplier_3_3`n + x_Supplier_3_4 + x_Supplier_3_5 - z_Supplier_3 >= 0`nMinAward_Supplier_4: x_Supplier_4_1 + x_Supplier_4_2 + x_Supplier_4_3`n + x_Supplier_4_4 + x_Supplier_4_5 - z_Supplier_4 >= 0`nMinAward_Supplier_5: x_Supplier_5_1 + x_Supplier_5_2 + x_Supplier_5_3`n + x_Supplier_5_4 + x_Supplier_5_5 - z_Supplier_5 >= 0`nOneDisc_Supplier_1: - z_Supplier_1 + zd_Supplier_1_0 = 0`nOneReb_Supplier_1: yr_Supplier_1_0 - z_Supplier_1 = 0`nRMax_Supplier_1_0: x_Supplier_1_1 + x_Supplier_1_2 + x_Supplier_1_3`n + x_Supplier_1_4 + x_Supplier_1_5 + 1000000000 yr_Supplier_1_0 <= 1001000000`nRMin_Supplier_1_0: x_Supplier_1_1 + x_Supplier_1_2 + x_Supplier_1_3`n + x_Supplier_1_4 + x_Supplier_1_5 - yr_Supplier_1_0 >= 0`nS0_Supplier_1: S0_Supplier_1 - 34 xKBX_Supplier_1_1 - 15 xKBX_Supplier_1_2`n - 15 xKBX_Supplier_1_3 - 75 xKBX_Supplier_1_4 - 24 xKBX_Supplier_1_5`n - 36 xSF_Supplier_1_1 - 17 xSF_Supplier_1_2 - 17 xSF_Supplier_1_3`n - 77 xSF_Supplier_1_4 - 26 xSF_Supplier_1_5 = 0`nS0_Supplier_2: S0_Supplier_2 - 14 xKBX_Supplier_2_1 - 15 xKBX_Supplier_2_2`n - 78 xKBX_Supplier_2_3 - 34 xKBX_Supplier_2_4 - 15 xKBX_Supplier_2_5`n - 16 xSF_Supplier_2_1 - 17 xSF_Supplier_2_2 - 80 xSF_Supplier_2_3`n - 36 xSF_Supplier_2_4 - 17 xSF_Supplier_2_5 = 0`nS0_Supplier_3: S0_Supplier_3 - 75 xKBX_Supplier_3_1 - 25 xKBX_Supplier_3_2`n - 56 xKBX_Supplier_3_3 - 24 xKBX_Supplier_3_4 - 87 xKBX_Supplier_3_5`n - 77 xSF_Supplier_3_1 - 27 xSF_Supplier_3_2 - 58 xSF_Supplier_3_3`n - 26 xSF_Supplier_3_4 - 89 xSF_Supplier_3_5 = 0`nS0_Supplier_4: S0_Supplier_4 - 93 xKBX_Supplier_4_1 - 24 xKBX_Supplier_4_2`n - 78 xKBX_Supplier_4_3 - 56 xKBX_Supplier_4_4 - 86 xKBX_Supplier_4_5`n - 95 xSF_Supplier_4_1 - 26 xSF_Supplier_4_2 - 80 xSF_Supplier_4_3`n - 58 xSF_Supplier_4_4 - 88 xSF_Supplier_4_5 = 0`nS0_Supplier_5: S0_Supplier_5 - 76 xKBX_Supplier_5_1 - 32 xKBX_Supplier_5_2`n - 89 xKBX_Supplier_5_3 - 13 xKBX_Supplier_5_4 - 68 xKBX_Supplier_5_5`n - 78 xSF_Supplier_5_1 - 34 xSF_Supplier_5_2 - 91 xSF_Supplier_5_3`n - 15 xSF_Supplier_5_4 - 70 xSF_Supplier_5_5 = 0`nSpend_Supplier_1: - S0_Supplier_1 + S_Supplier_1 + d_Supplier_1 = 0`nSpend_Supplier_2: - S0_Supplier_2 + S_Supplier_2 + d_Supplier_2 = 0`nSpend_Supplier_3: - S0_Supplier_3 + S_Supplier_3 + d_Supplier_3 = 0`nSpend_Supplier_4: - S0_Supplier_4 + S_Supplier_4 + d_Supplier_4 = 0`nSpend_Supplier_5: - S0_Supplier_5 + S_Supplier_5 + d_Supplier_5 = 0`nSupCntUB_0: Y_sup_0_Supplier_1 + Y_sup_0_Supplier_2 + Y_sup_0_Supplier_3`n + Y_sup_0_Supplier_4 + Y_sup_0_Supplier_5 <= 2`nTrLB_1_Supplier_1: x_Supplier_1_1 >= 0`nTrLB_1_Supplier_3: x_Supplier_3_1 >= 0`nTrLB_1_Supplier_4: x_Supplier_4_1 >= 0`nTrLB_1_Supplier_5: x_Supplier_5_1 >= 0`nTrLB_2_Supplier_1: x_Supplier_1_2 >= 0`nTrLB_2_Supplier_2: x_Supplier_2_2 >= 0`nTrLB_2_Supplier_4: x_Supplier_4_2 >= 0`nTrLB_2_Supplier_5: x_Supplier_5_2 >= 0`nTrLB_3_Supplier_1: x_Supplier_1_3 >= 0`nTrLB_3_Supplier_3: x_Supplier_3_3 >= 0`nTrLB_3_Supplier_4: x_Supplier_4_3 >= 0`nTrLB_3_Supplier_5: x_Supplier_5_3 >= 0`nTrLB_4_Supplier_2: x_Supplier_2_4 >= 0`nTrLB_4_Supplier_3: x_Supplier_3_4 >= 0`nTrLB_4_Supplier_4: x_Supplier_4_4 >= 0`nTrLB_4_Supplier_5: x_Supplier_5_4 >= 0`nTrLB_5_Supplier_1: x_Supplier_1_5 >= 0`nTrLB_5_Supplier_2: x_Supplier_2_5 >= 0`nTrLB_5_Supplier_4: x_Supplier_4_5 >= 0`nTrLB_5_Supplier_5: x_Supplier_5_5 >= 0`nTrUB_1_Supplier_1: - 1456 T_1_Supplier_1 + x_Supplier_1_1 <= 0`nTrUB_1_Supplier_3: - 1456 T_1_Supplier_3 + x_Supplier_3_1 <= 0`nTrUB_1_Supplier_4: - 1456 T_1_Supplier_4 + x_Supplier_4_1 <= 0`nTrUB_1_Supplier_5: - 1456 T_1_Supplier_5 + x_Supplier_5_1 <= 0`nTrUB_2_Supplier_1: - 2421 T_2_Supplier_1 + x_Supplier_1_2 <= 0`nTrUB_2_Supplier_2: - 2421 T_2_Supplier_2 + x_Supplier_2_2 <= 0`nTrUB_2_Supplier_4: - 2421 T_2_Supplier_4 + x_Supplier_4_2 <= 0`nTrUB_2_Supplier_5: - 2421 T_2_Supplier_5 + x_Supplier_5_2 <= 0`nTrUB_3_Supplier_1: - 6462 T_3_Supplier_1 + x_Supplier_1_3 <= 0`nTrUB_3_Supplier_3: - 6462 T_3_Supplier_3 + x_Supplier_3_3 <= 0`nTrUB_3_Supplier_4: - 6462 T_3_Supplier_4 + x_Supplier_4_3 <= 0`nTrUB_3_Supplier_5: - 6462 T_3_Supplier_5 + x_Supplier_5_3 <= 0`nTrUB_4_Supplier_2: - 2357 T_4_Supplier_2 + x_Supplier_2_4 <= 0`nTrUB_4_Supplier_3: - 2357 T_4_Supplier_3 + x_Supplier_3_4 <= 0`nTrUB_4_Supplier_4: - 2357 T_4_Supplier_4 + x_Supplier_4_4 <= 0`nTrUB_4_Supplier_5: - 2357 T_4_Supplier_5 + x_Supplier_5_4 <= 0`nTrUB_5_Supplier_1: - 2421 T_5_Supplier_1 + x_Supplier_1_5 <= 0`nTrUB_5_Supplier_2: - 2421 T_5_Supplier_2 + x_Supplier_2_5 <= 0`nTrUB_5_Supplier_4: - 2421 T_5_Supplier_4 + x_Supplier_4_5 <= 0`nTrUB_5_Supplier_5: - 2421 T_5_Supplier_5 + x_Supplier_5_5 <= 0`nYlb_0_Supplier_1: - 1e-06 Y_sup_0_Supplier_1 + x_Supplier_1_1 + x_Supplier_1_2`n + x_Supplier_1_3 + x_Supplier_1_4 + x_Supplier_1_5 >= 0`nYlb_0_Supplier_2: - 1e-06 Y_sup_0_Supplier_2 + x_Supplier_2_1 + x_Supplier_2_2`n + x_Supplier_2_3 + x_Supplier_2_4 + x_Supplier_2_5 >= 0`nYlb_0_Supplier_3: - 1e-06 Y_sup_0_Supplier_3 + x_Supplier_3_1 + x_Supplier_3_2`n + x_Supplier_3_3 + x_Supplier_3_4 + x_Supplier_3_5 >= 0`nYlb_0_Supplier_4: - 1e-06 Y_sup_0_Supplier_4 + x_Supplier_4_1 + x_Supplier_4_2`n + x_Supplier_4_3 + x_Supplier_4_4 + x_Supplier_4_5 >= 0`nYlb_0_Supplier_5: - 1e-06 Y_sup_0_Supplier_5 + x_Supplier_5_1 + x_Supplier_5_2`n + x_Supplier_5_3 + x_Supplier_5_4 + x_Supplier_5_5 >= 0`nYub_0_Supplier_1: - 1000000000 Y_sup_0_Supplier_1 + x_Supplier_1_1`n + x_Supplier_1_2 + x_Supplier_1_3 + x_Supplier_1_4 + x_Supplier_1_5 <= 0`nYub_0_Supplier_2: - 1000000000 Y_sup_0_Supplier_2 + x_Supplier_2_1`n + x_Supplier_2_2 + x_Supplier_2_3 + x_Supplier_2_4 + x_Supplier_2_5 <= 0`nYub_0_Supplier_3: - 1000000000 Y_sup_0_Supplier_3 + x_Supplier_3_1`n + x_Supplier_3_2 + x_Supplier_3_3 + x_Supplier_3_4 + x_Supplier_3_5 <= 0`nYub_0_Supplier_4: - 1000000000 Y_sup_0_Supplier_4 + x_Supplier_4_1`n + x_Supplier_4_2 + x_Supplier_4_3 + x_Supplier_4_4 + x_Supplier_4_5 <= 0`nYub_0_Supplier_5: - 1000000000 Y_sup_0_Supplier_5 + x_Supplier_5_1`n + x_Supplier_5_2 + x_Supplier_5_3 + x_Supplier_5_4 + x_Supplier_5_5 <= 0`ndLow_Supplier_1_0: - 0.02 S0_Supplier_1 + d_Supplier_1`n - 1406346 zd_Supplier_1_0 >= -1406346`ndUp_Supplier_1_0: - 0.02 S0_Supplier_1 + d_Supplier_1`n + 1406346 zd_Supplier_1_0 <= 1406346`ndZeroNoTier_Supplier_2: d_Supplier_2 = 0`ndZeroNoTier_Supplier_3: d_Supplier_3 = 0`ndZeroNoTier_Supplier_4: d_Supplier_4 = 0`ndZeroNoTier_Supplier_5: d_Supplier_5 = 0`ndZero_Supplier_1: d_Supplier_1 - 1406346 zd_Supplier_1_0 <= 0`nrLow_Supplier_1_0: - 0.95 S_Supplier_1 + reb_Supplier_1`n - 1406346 yr_Supplier_1_0 >= -1406346`nrUp_Supplier_1_0: - 0.95 S_Supplier_1 + reb_Supplier_1`n + 1406346 yr_Supplier_1_0 <= 1406346`nrZeroNoTier_Supplier_2: reb_Supplier_2 = 0`nrZeroNoTier_Supplier_3: reb_Supplier_3 = 0`nrZeroNoTier_Supplier_4: reb_Supplier_4 = 0`nrZeroNoTier_Supplier_5: reb_Supplier_5 = 0`nrZero_Supplier_1: reb_Supplier_1 - 1406346 yr_Supplier_1_0 <= 0`nBinaries`nT_1_Supplier_1`nT_1_Supplier_3`nT_1_Supplier_4`nT_1_Supplier_5`nT_2_Supplier_1`nT_2_Supplier_2`nT_2_Supplier_4`nT_2_Supplier_5`nT_3_Supplier_1`nT_3_Supplier_3`nT_3_Supplier_4`nT_3_Supplier_5`nT_4_Supplier_2`nT_4_Supplier_3`nT_4_Supplier_4`nT_4_Supplier_5`nT_5_Supplier_1`nT_5_Supplier_2`nT_5_Supplier_4`nT_5_Supplier_5`nY_sup_0_Supplier_1`nY_sup_0_Supplier_2`nY_sup_0_Supplier_3`nY_sup_0_Supplier_4`nY_sup_0_Supplier_5`nbDDP_Supplier_1_1`nbDDP_Supplier_1_2`nbDDP_Supplier_1_3`nbDDP_Supplier_1_4`nbDDP_Supplier_1_5`nbDDP_Supplier_2_1`nbDDP_Supplier_2_2`nbDDP_Supplier_2_3`nbDDP_Supplier_2_4`nbDDP_Supplier_2_5`nbDDP_Supplier_3_1`nbDDP_Supplier_3_2`nbDDP_Supplier_3_3`nbDDP_Supplier_3_4`nbDDP_Supplier_3_5`nbDDP_Supplier_4_1`nbDDP_Supplier_4_2`nbDDP_Supplier_4_3`nbDDP_Supplier_4_4`nbDDP_Supplier_4_5`nbDDP_Supplier_5_1`nbDDP_Supplier_5_2`nbDDP_Supplier_5_3`nbDDP_Supplier_5_4`nbDDP_Supplier_5_5`nyr_Supplier_1_0`nz_Supplier_1`nz_Supplier_2`nz_Supplier_3`nz_Supplier_4`nz_Supplier_5`nzd_Supplier_1_0`nEnd`n"
